# Insert a new weekly price-report row for "Feria Lagunitas de Puerto Montt" / Repollo.
# The new row is inserted at row 244, pushing the former rows 244-280 down to 245-281
# (dimension grows from A1:R280 to A1:R281). The new row carries the same
# market/category/quality metadata as the (now-shifted) row below it, but with its own
# date and volume values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(244).Insert()

$ws.Range("A244").Value = 4
$ws.Range("B244").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C244").Value = 'Los Lagos'
$ws.Range("D244").Value = 44505
$ws.Range("E244").Value = 10
$ws.Range("F244").Value = 100112006
$ws.Range("G244").Value = 'Repollo'
$ws.Range("H244").Value = 'Crespo record'
$ws.Range("I244").Value = 'Primera'
$ws.Range("J244").Value = 500
$ws.Range("K244").Value = 1200
$ws.Range("L244").Value = 1200
$ws.Range("M244").Value = 1200
$ws.Range("N244").Value = '$/unidad'
$ws.Range("O244").Value = 'Región de Coquimbo'
$ws.Range("P244").Value = 1200
$ws.Range("Q244").Value = 1
$ws.Range("R244").Value = 'Hortaliza'
